# books.xlsx refactor:
#  - de-duplicate the doubled author-name strings
#  - reset the header/body cell alignment back to Excel defaults
#    (general / bottom / no wrap)
#  - shrink column F and drop the custom row heights
#  - remove the trailing blank row (row 8)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Fix the duplicated "Author name(s)" text in column F (shared strings
#    had the author text repeated twice back-to-back).
# ---------------------------------------------------------------------
$ws.Range("F2").Value = "Ronald L. Graham; Donald Knuth; Oren Patashnik"
$ws.Range("F3").Value = "Ronald L. Graham; Donald Knuth; Oren Patashnik"
$ws.Range("F4").Value = "Gayle Laakmann McDowell"
$ws.Range("F5").Value = "Gayle Laakmann McDowell"
$ws.Range("F6").Value = "Albert Einstein"
$ws.Range("F7").Value = "Charles Darwin"

# ---------------------------------------------------------------------
# 2. Reset alignment on the used range back to Excel's defaults:
#    horizontal=general, vertical=bottom, wrapText=off.
# ---------------------------------------------------------------------
$rng = $ws.Range("A1:H7")
$rng.HorizontalAlignment = 1
$rng.VerticalAlignment = -4107
$rng.WrapText = $False

# ---------------------------------------------------------------------
# 3. Narrow column F and drop the explicit 30pt row heights (rows 1-7)
#    back to the sheet default.
# ---------------------------------------------------------------------
$ws.Columns("F").ColumnWidth = 44.67
$ws.Rows("1:7").AutoFit()

# ---------------------------------------------------------------------
# 4. Remove the trailing empty row 8.
# ---------------------------------------------------------------------
$ws.Rows("8").Delete()
